$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Comercializadora del Agro de
# Limarí" (Haba). It belongs chronologically right after the existing
# row 26, so insert a fresh row there and push everything else down by
# one (rows 26-64 -> 27-65).
$ws.Rows.Item(26).Insert()

# Fill in the new record's data (same market/category/quality metadata
# as the rest of the block; only the date, volume and price columns are
# specific to this record).
$ws.Cells.Item(26, 1).Value = 2
$ws.Cells.Item(26, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44790
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100112026
$ws.Cells.Item(26, 7).Value = "Haba"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 1460
$ws.Cells.Item(26, 11).Value = 9000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 13).Value = 9500
$ws.Cells.Item(26, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 380
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
